$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Title cell (A1): rename the generic template title to the real project title.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value2 = "Mood Badge -- V1 / Megan Singer"

# ---------------------------------------------------------------------------
# New BOM rows (10-39) listing every part name / designator used in the design.
# Column order: B=Comment, C=Designator, E=Mounting Type, F=Package, G=Manufacturer,
# H=Manufacturer Part Number
# ---------------------------------------------------------------------------

# Row 10 - 100 nF capacitor
$ws.Range("B10").Value2 = "100 nF capacitor"
$ws.Range("C10").Value2 = "C2, C3, C4, C5, C6, C7, C15, C21"
$ws.Range("E10").Value2 = "SMD"
$ws.Range("F3").Copy()
$ws.Range("F10").PasteSpecial(-4122)
$ws.Range("F10").Value2 = "'0603"

# Row 11 - 4.7 uF capacitor
$ws.Range("B11").Value2 = "4.7 uF capacitor"
$ws.Range("C11").Value2 = "C1, C8, C9, C16, C17, C18, "
$ws.Range("E11").Value2 = "SMD"
$ws.Range("F3").Copy()
$ws.Range("F11").PasteSpecial(-4122)
$ws.Range("F11").Value2 = "'0603"

# Row 12 - 2.2 uF capacitor
$ws.Range("B12").Value2 = "2.2 uF capacitor"
$ws.Range("C12").Value2 = "C10, C19"
$ws.Range("E12").Value2 = "SMD"
$ws.Range("F3").Copy()
$ws.Range("F12").PasteSpecial(-4122)
$ws.Range("F12").Value2 = "'0603"

# Row 13 - 1 uF capacitor
$ws.Range("B13").Value2 = "1 uF capacitor"
$ws.Range("C13").Value2 = "C20"
$ws.Range("E13").Value2 = "SMD"
$ws.Range("F3").Copy()
$ws.Range("F13").PasteSpecial(-4122)
$ws.Range("F13").Value2 = "'0603"

# Row 14 - 0.1 uF capacitor
$ws.Range("B14").Value2 = "0.1 uF capacitor"
$ws.Range("C14").Value2 = "C11, C14"
$ws.Range("E14").Value2 = "SMD"
$ws.Range("F3").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("F14").Value2 = "'0603"

# Row 15 - 22 pF capacitor
$ws.Range("B15").Value2 = "22 pF capacitor"
$ws.Range("C15").Value2 = "C12, C13"
$ws.Range("E15").Value2 = "SMD"
$ws.Range("F3").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("F15").Value2 = "'0603"

# Row 16 - 205 ohm resistor
$ws.Range("B16").Value2 = "205 ohm resistor"
$ws.Range("C16").Value2 = "R16, R17, R18, R19, R20, R21, R22, R23"
$ws.Range("E16").Value2 = "SMD"
$ws.Range("F3").Copy()
$ws.Range("F16").PasteSpecial(-4122)
$ws.Range("F16").Value2 = "'0603"

# Row 17 - 487 ohm resistor
$ws.Range("B17").Value2 = "487 ohm resistor"
$ws.Range("C17").Value2 = "R24"
$ws.Range("E17").Value2 = "SMD"
$ws.Range("F3").Copy()
$ws.Range("F17").PasteSpecial(-4122)
$ws.Range("F17").Value2 = "'0603"

# Row 18 - 1 kohm resistor
$ws.Range("B18").Value2 = "1 kohm resistor"
$ws.Range("C18").Value2 = "R2, R4, R6, R8, R10, R12, R13"
$ws.Range("E18").Value2 = "SMD"
$ws.Range("F3").Copy()
$ws.Range("F18").PasteSpecial(-4122)
$ws.Range("F18").Value2 = "'0603"

# Row 19 - 10 kohm resistor
$ws.Range("B19").Value2 = "10 kohm resistor"
$ws.Range("C19").Value2 = "R1, R3, R5, R7, R9, R11, R14, R15, R26"
$ws.Range("E19").Value2 = "SMD"
$ws.Range("F3").Copy()
$ws.Range("F19").PasteSpecial(-4122)
$ws.Range("F19").Value2 = "'0603"

# Row 20 - 12 kohm resistor
$ws.Range("B20").Value2 = "12 kohm resistor"
$ws.Range("C20").Value2 = "R25"
$ws.Range("E20").Value2 = "SMD"
$ws.Range("F3").Copy()
$ws.Range("F20").PasteSpecial(-4122)
$ws.Range("F20").Value2 = "'0603"

# Row 21 - 210 kohm resistor
$ws.Range("B21").Value2 = "210 kohm resistor"
$ws.Range("C21").Value2 = "R28, R29"
$ws.Range("E21").Value2 = "SMD"
$ws.Range("F3").Copy()
$ws.Range("F21").PasteSpecial(-4122)
$ws.Range("F21").Value2 = "'0603"

# Row 22 - 3.3 Mohm resistor
$ws.Range("B22").Value2 = "3.3 Mohm resistor"
$ws.Range("C22").Value2 = "R27"
$ws.Range("E22").Value2 = "SMD"
$ws.Range("F3").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("F22").Value2 = "'0603"

# Row 23 - Clock crystal
$ws.Range("B23").Value2 = "Clock crystal"
$ws.Range("C23").Value2 = "Y1"
$ws.Range("E23").Value2 = "THT"
$ws.Range("F3").Copy()
$ws.Range("F23").PasteSpecial(-4122)
$ws.Range("H23").Value2 = "CA-310"

# Row 24 - Diode
$ws.Range("B24").Value2 = "Diode"
$ws.Range("C24").Value2 = "D1, D2, D3"
$ws.Range("E24").Value2 = "SMD"
$ws.Range("F3").Copy()
$ws.Range("F24").PasteSpecial(-4122)
$ws.Range("F24").Value2 = "'0603"

# Row 25 - IR LED
$ws.Range("B25").Value2 = "IR LED"
$ws.Range("C25").Value2 = "LED8"
$ws.Range("F3").Copy()
$ws.Range("F25").PasteSpecial(-4122)

# Row 26 - LED
$ws.Range("B26").Value2 = "LED"
$ws.Range("C26").Value2 = "LED1, LED2, LED3, LED4, LED5, LED6, LED7"
$ws.Range("E26").Value2 = "SMD"
$ws.Range("F3").Copy()
$ws.Range("F26").PasteSpecial(-4122)
$ws.Range("F26").Value2 = "'0603"

# Row 27 - LED matrix
$ws.Range("B27").Value2 = "LED matrix"
$ws.Range("C27").Value2 = "MATRIX"
$ws.Range("E27").Value2 = "THT"
$ws.Range("F3").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("G27").Value2 = "Adafruit Industries"
$ws.Range("H27").Value2 = "SEGMENT_8X8_ROWCATHODEBL-M12A883"

# Row 28 - Large push button
$ws.Range("B28").Value2 = "Large push button"
$ws.Range("C28").Value2 = "UP, DOWN, LEFT, RIGHT, SELECT, MODE"
$ws.Range("E28").Value2 = "SMD"
$ws.Range("H28").Value2 = "SPST.KT11P4SM34LFS"

# Row 29 - Small push button
$ws.Range("B29").Value2 = "Small push button"
$ws.Range("C29").Value2 = "RESET"
$ws.Range("E29").Value2 = "SMD"
$ws.Range("H29").Value2 = "SPST-PTS830GG140"

# Row 30 - Switch
$ws.Range("B30").Value2 = "Switch"
$ws.Range("C30").Value2 = "ON/OFF"
$ws.Range("E30").Value2 = "SMD"
$ws.Range("H30").Value2 = "SWITCH.JS202011JXQNUPWARDS"

# Row 31 - Bootloader 2x3 pins
$ws.Range("B31").Value2 = "Bootloader 2x3 pins"
$ws.Range("C31").Value2 = "BOOT"
$ws.Range("E31").Value2 = "THT"

# Row 32 - FTDI 1x6 pins
$ws.Range("B32").Value2 = "FTDI 1x6 pins"
$ws.Range("C32").Value2 = "FTDI_DEVICE"
$ws.Range("E32").Value2 = "THT"

# Row 33 - USB port
$ws.Range("B33").Value2 = "USB port"
$ws.Range("C33").Value2 = "USB1"
$ws.Range("E33").Value2 = "SMD"

# Row 34 - LIPO battery port
$ws.Range("B34").Value2 = "LIPO battery port"
$ws.Range("C34").Value2 = "LIPO"
$ws.Range("E34").Value2 = "SMD"

# Row 35 - Battery charger IC
$ws.Range("B35").Value2 = "Battery charger IC"
$ws.Range("C35").Value2 = "CHRG"
$ws.Range("E35").Value2 = "SMD"
$ws.Range("H35").Value2 = "MCP73831"

# Row 36 - Port expander
$ws.Range("B36").Value2 = "Port expander"
$ws.Range("C36").Value2 = "IC1"
$ws.Range("E36").Value2 = "SMD"
$ws.Range("F36").Value2 = "SO16W"
$ws.Range("H36").Value2 = "PCF8574T"

# Row 37 - Arduino ATMEGA328
$ws.Range("B37").Value2 = "Arduino ATMEGA328"
$ws.Range("C37").Value2 = "U2"
$ws.Range("E37").Value2 = "SMD"
$ws.Range("F37").Value2 = "QFP"

# Row 38 - Op amp
$ws.Range("B38").Value2 = "Op amp"
$ws.Range("C38").Value2 = "OP"
$ws.Range("E38").Value2 = "SMD"
$ws.Range("F38").Value2 = "SOT23-5"

# Row 39 - Transducer
$ws.Range("B39").Value2 = "Transducer"
$ws.Range("C39").Value2 = "PD"
$ws.Range("E39").Value2 = "SMD"
$ws.Range("F39").Value2 = "APDS9008"

# ---------------------------------------------------------------------------
# Update the running-total formulas so they cover the newly added rows.
# ---------------------------------------------------------------------------
$ws.Range("M1").Formula = "=SUM(M3:M81)"
$ws.Range("O1").Formula = "=SUM(O3:O81)"

# Leave the selection where the author last left it.
$ws.Range("F24").Select()
